# Correction in SA algorithm and 746 logs
# Updates the "Fitness" column (C) values in the run_10 log sheet to reflect
# the corrected simulated-annealing results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-49 (Generation 0-47) -> 8077
$ws.Range("C2:C49").Value = 8077

# Rows 50-166 (Generation 48-164) -> 7622
$ws.Range("C50:C166").Value = 7622

# Rows 167-172 (Generation 165-170) -> 7295
$ws.Range("C167:C172").Value = 7295
